# Add a new "TestScripts" worksheet after the last existing sheet (Cookie),
# populate it with the TestName/Enabled header + one data row, size column A,
# and leave it as the active/selected sheet - matching the commit that added
# support for dynamically enabling/disabling test scripts via an
# IAnnotationTransformer.

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$testScripts = $wb.Worksheets.Add($null, $lastSheet)
$testScripts.Name = "TestScripts"

$testScripts.Range("A1").Value = "TestName"
$testScripts.Range("B1").Value = "Enabled"
$testScripts.Range("A2").Value = "AnnotationTtest "
$testScripts.Range("B2").Value = $true

# Column A width of 32 (COM's ColumnWidth adds ~0.8333 padding internally, so
# back that off to land on an on-disk width of exactly 32).
$testScripts.Columns.Item(1).ColumnWidth = 31.166666666666668

[void]$testScripts.Range("B2").Select()
[void]$testScripts.Activate()

# Sheet2 was the previously tab-selected sheet; now that TestScripts is
# active/selected that flag should no longer be set there (handled
# automatically by activating the new sheet above).

Write-Output "done"
